$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4032.7273
$ws.Range("I112").Value = 696.6667
$ws.Range("J112").Value = 4276.829
$ws.Range("K112").Value = 2090.0001
$ws.Range("L112").Value = 12830.487
$ws.Range("M112").Value = -982.0001000000002
$ws.Range("N112").Value = -15046.487

$ws.Range("H116").Value = 915655.9399999999
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 1118023.9
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 1118023.9
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -1124907.9

$ws.Range("H121").Value = 1200.9584
$ws.Range("J121").Value = 1187.9565
$ws.Range("L121").Value = 3563.8695
$ws.Range("N121").Value = -7057.8695

$ws.Range("H141").Value = 13904.444
$ws.Range("I141").Value = 15020
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 45060
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -39880
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4327.9707
$ws.Range("I61").Value = 5263.2383
$ws.Range("J61").Value = 2817.1538
$ws.Range("K61").Value = 5263.2383
$ws.Range("L61").Value = 2817.1538
$ws.Range("M61").Value = -5051.2383
$ws.Range("N61").Value = -3241.1538

$ws.Range("H74").Value = 1212.6923
$ws.Range("I74").Value = 802.5
$ws.Range("J74").Value = 1287.2727
$ws.Range("K74").Value = 802.5
$ws.Range("L74").Value = 1287.2727
$ws.Range("M74").Value = 71.5
$ws.Range("N74").Value = -3035.2727

$ws.Range("H77").Value = 1212.6923
$ws.Range("I77").Value = 802.5
$ws.Range("J77").Value = 1287.2727
$ws.Range("K77").Value = 4012.5
$ws.Range("L77").Value = 6436.363499999999
$ws.Range("M77").Value = 355.5
$ws.Range("N77").Value = -15172.3635

$ws.Range("H132").Value = 2978495.2
$ws.Range("I132").Value = 6580529
$ws.Range("K132").Value = 19741587
$ws.Range("M132").Value = -19739057

$ws.Range("H133").Value = 8882.625
$ws.Range("J133").Value = 8882.625
$ws.Range("L133").Value = 8882.625
$ws.Range("N133").Value = -13942.625

$ws.Range("H136").Value = 4327.9707
$ws.Range("I136").Value = 5263.2383
$ws.Range("J136").Value = 2817.1538
$ws.Range("K136").Value = 15789.7149
$ws.Range("L136").Value = 8451.4614
$ws.Range("M136").Value = -13239.7149
$ws.Range("N136").Value = -13551.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 747.2381
$ws.Range("I94").Value = 552.5
$ws.Range("J94").Value = 1234.0834
$ws.Range("K94").Value = 552.5
$ws.Range("L94").Value = 1234.0834
$ws.Range("M94").Value = -101.5
$ws.Range("N94").Value = -2136.0834

$ws.Range("H134").Value = 10117686
$ws.Range("I134").Value = 13910992
$ws.Range("J134").Value = 2201.5557
$ws.Range("K134").Value = 41732976
$ws.Range("L134").Value = 6604.6671
$ws.Range("M134").Value = -41730441
$ws.Range("N134").Value = -11674.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1481.8182
$ws.Range("K16").Value = 1481.8182
$ws.Range("M16").Value = -1194.8182

$ws.Range("H86").Value = 3095
$ws.Range("I86").Value = 3075
$ws.Range("J86").Value = 3115
$ws.Range("K86").Value = 3075
$ws.Range("L86").Value = 3115
$ws.Range("M86").Value = -1952
$ws.Range("N86").Value = -5361

$ws.Range("H89").Value = 3095
$ws.Range("I89").Value = 3075
$ws.Range("J89").Value = 3115
$ws.Range("K89").Value = 15375
$ws.Range("L89").Value = 15575
$ws.Range("M89").Value = -9759
$ws.Range("N89").Value = -26807

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1481.8182
$ws.Range("K113").Value = 1481.8182
$ws.Range("M113").Value = 688.1818000000001

$ws.Range("H132").Value = 6948111.5
$ws.Range("I132").Value = 12821532
$ws.Range("J132").Value = 6795.9546
$ws.Range("K132").Value = 38464596
$ws.Range("L132").Value = 20387.8638
$ws.Range("M132").Value = -38462066
$ws.Range("N132").Value = -25447.8638

$ws.Range("H133").Value = 19399.6
$ws.Range("J133").Value = 19399.6
$ws.Range("L133").Value = 19399.6
$ws.Range("N133").Value = -24459.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2124.7058
$ws.Range("I140").Value = 2048.889
$ws.Range("J140").Value = 2210
$ws.Range("K140").Value = 6146.667
$ws.Range("L140").Value = 6630
$ws.Range("M140").Value = -966.6670000000004
$ws.Range("N140").Value = -16990

$ws.Range("H141").Value = 3093.8635
$ws.Range("I141").Value = 1895.2778
$ws.Range("J141").Value = 8487.5
$ws.Range("K141").Value = 5685.8334
$ws.Range("L141").Value = 25462.5
$ws.Range("M141").Value = -505.8334000000004
$ws.Range("N141").Value = -35822.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27031516
$ws.Range("I132").Value = 58825710
$ws.Range("J132").Value = 6449.75
$ws.Range("K132").Value = 176477130
$ws.Range("L132").Value = 19349.25
$ws.Range("M132").Value = -176474600
$ws.Range("N132").Value = -24409.25

$ws.Range("H133").Value = 40750
$ws.Range("J133").Value = 40750
$ws.Range("L133").Value = 40750
$ws.Range("N133").Value = -50870

$ws.Range("H135").Value = 38428.57
$ws.Range("J135").Value = 38428.57
$ws.Range("L135").Value = 38428.57
$ws.Range("N135").Value = -48568.57

$ws.Range("H137").Value = 17400
$ws.Range("J137").Value = 17400
$ws.Range("L137").Value = 17400
$ws.Range("N137").Value = -27600

$ws.Range("H138").Value = 61331.668
$ws.Range("J138").Value = 61331.668
$ws.Range("L138").Value = 61331.668
$ws.Range("N138").Value = -71611.66800000001

$ws.Range("H140").Value = 15796.667
$ws.Range("J140").Value = 15796.667
$ws.Range("L140").Value = 15796.667
$ws.Range("N140").Value = -26156.667

$ws.Range("H141").Value = 11923.7
$ws.Range("J141").Value = 11923.7
$ws.Range("L141").Value = 11923.7
$ws.Range("N141").Value = -22283.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1549.9546
$ws.Range("I7").Value = 1509.95
$ws.Range("J7").Value = 1950
$ws.Range("K7").Value = 1509.95
$ws.Range("L7").Value = 1950
$ws.Range("M7").Value = -1397.95
$ws.Range("N7").Value = -2174

$ws.Range("H126").Value = 1549.9546
$ws.Range("I126").Value = 1509.95
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 4529.85
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -2059.85
$ws.Range("N126").Value = -10790

$ws.Range("H132").Value = 8698861
$ws.Range("I132").Value = 18184518
$ws.Range("J132").Value = 3674
$ws.Range("K132").Value = 54553554
$ws.Range("L132").Value = 11022
$ws.Range("M132").Value = -54551024
$ws.Range("N132").Value = -16082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 3938
$ws.Range("J41").Value = 5377
$ws.Range("L41").Value = 5377
$ws.Range("N41").Value = -6157

$ws.Range("H45").Value = 4807.125
$ws.Range("J45").Value = 4807.125
$ws.Range("L45").Value = 4807.125
$ws.Range("N45").Value = -5789.125

$ws.Range("H132").Value = 24655240
$ws.Range("I132").Value = 18280480
$ws.Range("J132").Value = 28550928
$ws.Range("K132").Value = 54841440
$ws.Range("L132").Value = 85652784
$ws.Range("M132").Value = -54838910
$ws.Range("N132").Value = -85657844
